$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44162
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 5000
$ws.Range("P2").Value = 417

# Row 3
$ws.Range("D3").Value = 44329
$ws.Range("J3").Value = 40
$ws.Range("K3").Value = 5000
$ws.Range("L3").Value = 6000
$ws.Range("M3").Value = 5500
$ws.Range("P3").Value = 458

# Row 5
$ws.Range("D5").Value = 44455
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 7000
$ws.Range("P5").Value = 583

# Row 6
$ws.Range("D6").Value = 44442
$ws.Range("J6").Value = 20

# Row 7
$ws.Range("D7").Value = 44441
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 6000
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = 6000
$ws.Range("P7").Value = 500

# Row 8
$ws.Range("D8").Value = 44179
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 6000
$ws.Range("P8").Value = 500

# Row 9
$ws.Range("D9").Value = 44302
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 6000
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = 6000
$ws.Range("P9").Value = 500

# Row 10
$ws.Range("D10").Value = 44196
$ws.Range("J10").Value = 20
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = 5000
$ws.Range("P10").Value = 417

# Row 11
$ws.Range("D11").Value = 44299
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("P11").Value = 583

# Row 12
$ws.Range("D12").Value = 44211
$ws.Range("J12").Value = 65
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = 5000
$ws.Range("P12").Value = 417

# Row 13
$ws.Range("D13").Value = 44327
$ws.Range("K13").Value = 6000
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = 6000
$ws.Range("P13").Value = 500

# Row 14
$ws.Range("D14").Value = 44424
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 7000
$ws.Range("M14").Value = 7000
$ws.Range("P14").Value = 583

# Row 16
$ws.Range("D16").Value = 44428
$ws.Range("J16").Value = 10

# Row 17
$ws.Range("D17").Value = 44195
$ws.Range("J17").Value = 55
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("M17").Value = 5000
$ws.Range("P17").Value = 417

# Row 18
$ws.Range("D18").Value = 44454
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 7000
$ws.Range("P18").Value = 583

# Row 19
$ws.Range("D19").Value = 44410

# Row 20
$ws.Range("D20").Value = 44467
$ws.Range("J20").Value = 20

# Row 21
$ws.Range("D21").Value = 44427
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 7000
$ws.Range("P21").Value = 583

# Row 22
$ws.Range("D22").Value = 44411

# Row 23
$ws.Range("D23").Value = 44203
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("M23").Value = 5000
$ws.Range("P23").Value = 417

# Row 24
$ws.Range("D24").Value = 44413
$ws.Range("J24").Value = 40

# Row 26
$ws.Range("D26").Value = 44474

# Row 27
$ws.Range("D27").Value = 44186
$ws.Range("K27").Value = 5000
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = 5000
$ws.Range("P27").Value = 417

# Row 28
$ws.Range("D28").Value = 44372
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 6000
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = 6000
$ws.Range("P28").Value = 500

# Row 29
$ws.Range("D29").Value = 44448
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 7000
$ws.Range("M29").Value = 7000
$ws.Range("P29").Value = 583

# Row 30
$ws.Range("D30").Value = 44452
$ws.Range("J30").Value = 40
$ws.Range("K30").Value = 7000
$ws.Range("L30").Value = 7000
$ws.Range("M30").Value = 7000
$ws.Range("P30").Value = 583

# Row 31
$ws.Range("D31").Value = 44453
$ws.Range("K31").Value = 7000
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = 7000
$ws.Range("O31").Value = "Provincia de Cautín"
$ws.Range("P31").Value = 583

# Row 32
$ws.Range("D32").Value = 44435
$ws.Range("J32").Value = 30

# Row 33
$ws.Range("D33").Value = 44166
$ws.Range("J33").Value = 55

# Row 34
$ws.Range("D34").Value = 44369
$ws.Range("J34").Value = 20
$ws.Range("K34").Value = 4000
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = 4000
$ws.Range("O34").Value = "Región Metropolitana"
$ws.Range("P34").Value = 333

# Row 35
$ws.Range("D35").Value = 44301
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 6000
$ws.Range("L35").Value = 6000
$ws.Range("M35").Value = 6000
$ws.Range("P35").Value = 500

# Row 36
$ws.Range("D36").Value = 44425
$ws.Range("J36").Value = 20
$ws.Range("K36").Value = 7000
$ws.Range("L36").Value = 7000
$ws.Range("M36").Value = 7000
$ws.Range("P36").Value = 583

# Row 37
$ws.Range("D37").Value = 44466
$ws.Range("J37").Value = 50
$ws.Range("K37").Value = 6000
$ws.Range("M37").Value = 6400
$ws.Range("P37").Value = 533

# Row 38
$ws.Range("D38").Value = 44438
$ws.Range("J38").Value = 30
$ws.Range("K38").Value = 6000
$ws.Range("L38").Value = 6000
$ws.Range("M38").Value = 6000
$ws.Range("P38").Value = 500

# Row 39
$ws.Range("D39").Value = 44326
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 6000
$ws.Range("M39").Value = 6000
$ws.Range("P39").Value = 500

# Row 40
$ws.Range("D40").Value = 44432
$ws.Range("J40").Value = 30
